# Re-sort / arrange the item list alphabetically by brand and add the new
# "Sk-Mox" brand + item (DB connections updated -> new item pulled in,
# whole sheet re-sorted/arranged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the two Etorix sub-items (rows 7 & 8), and slot the third one in ---
$ws.Range("D7").Value = "Etorix 90mg Tablet"
$ws.Range("E7").Value = "30's"

$ws.Range("D8").Value = "Etorix 120mg Tablet"
$ws.Range("E8").Value = "20's"

$ws.Range("D9").Value = "Etorix 60mg Tablet - 40's"
$ws.Range("E9").Value = "40's"

# --- Swap the two Flucloxin sub-items (rows 11 & 12) ---
$ws.Range("D11").Value = "Flucloxin 500mg Capsule - 36's"
$ws.Range("E11").Value = "36 's"

$ws.Range("D12").Value = "Flucloxin 500mg Capsule"
$ws.Range("E12").Value = "30 's"

# --- Re-arrange the three Ketonic sub-items (rows 14, 15, 16) ---
$ws.Range("D14").Value = "Ketonic 10mg Tablet"
$ws.Range("E14").Value = "20's"

$ws.Range("D15").Value = "Ketonic 30mg Injection"
$ws.Range("E15").Value = "5 's"

$ws.Range("D16").Value = "Ketonic 30mg IM/IV Injection - 4's"
$ws.Range("E16").Value = "4's"

# --- Swap the two Kynol sub-items (rows 17 & 18) ---
$ws.Range("D17").Value = "Kynol D 25mg Tablet"
$ws.Range("E17").Value = "60 's"

$ws.Range("D18").Value = "Kynol TR 200mg Capsule"
$ws.Range("E18").Value = "30 's"

# --- Row 24 becomes the new "Sk-Mox" brand/item (inserted via the DB refresh) ---
$ws.Range("A24").Value = 24
$ws.Range("B24").Value = "Sk-Mox"
$ws.Range("D24").Value = "Sk-Mox 500mg Capsule"
$ws.Range("E24").Value = "48 's"

# --- Zithrox items shift down one row to make room, row 25 now holds what
#     used to be in row 24's "15ml Suspension" slot ---
$ws.Range("D25").Value = "Zithrox 15ml Suspension"
$ws.Range("E25").Value = "15 ml"

# --- New row 28: the former row-25 "30ml Dry Suspension" item, now at the end ---
$ws.Range("A28").Value = 35
$ws.Range("B28").Value = "Zithrox"
$ws.Range("C28").Value = 27
$ws.Range("D28").Value = "Zithrox 30ml Dry Suspension"
$ws.Range("E28").Value = "30ml"
